# Update "想去人数" (people interested count) values in the "展览" and
# "全部类型" sheets to reflect the latest scrape, as described in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    2  = 234
    4  = 527
    5  = 13886
    8  = 8
    9  = 1784
    10 = 171
    14 = 535
    17 = 13939
    19 = 625
    20 = 14967
    22 = 8264
    25 = 25
    28 = 165
    31 = 1038
    32 = 13
    33 = 20
    39 = 217
    42 = 5088
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoUpdates[$row]
}

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 234
    4  = 527
    5  = 13886
    8  = 8
    9  = 1784
    10 = 171
    14 = 535
    17 = 13939
    18 = 368
    19 = 625
    20 = 14967
    22 = 8264
    25 = 25
    28 = 165
    31 = 1038
    32 = 14
    33 = 20
    41 = 217
    44 = 5088
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
